$d = $word.ActiveDocument

# Locate the paragraph ("11/4/16 finished accurate turn ...") that needs:
#   1. a paragraph-mark run-properties block (sz/szCs 32) added to its pPr
#   2. a brand new paragraph inserted right after it containing the
#      11/11/16 follow-up note (the _GoBack bookmark travels with the new,
#      last paragraph, same as it would after a real edit).
$p3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*11/4/16 finished accurate turn*") {
        $p3 = $candidate
    }
}
if ($p3 -eq $null) {
    $p3 = $d.Paragraphs.Item(3)
}

# Target the paragraph's range but drop the trailing paragraph mark so the
# replacement XML supplies its own paragraph boundaries cleanly.
$r = $d.Range($p3.Range.Start, $p3.Range.End - 1)

$newXml = '<?xml version="1.0" standalone="yes"?>' +
  '<?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
              '<w:r><w:t xml:space="preserve">11/4/16 finished accurate turn my block, need to add comments. Has been tested, lives up to name. </w:t></w:r>' +
              '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:sym w:font="Wingdings" w:char="F04A"/></w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:r><w:t xml:space="preserve">11/11/16 fixed accurate turn bug just had to reset timer is done variable. Bug was first accurate turn worked second skipped through turn part because the timer is done variable automatically ended loop. Now it truly lives up to name.  </w:t></w:r>' +
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$r.InsertXML($newXml) | Out-Null
